$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price ("Price", column D) and
# volume-change ("Volume(1h)", column E) figures scraped for this run.
# Values are written as literal text (matching the original inline-string
# cells) by briefly switching the cell to Text format while assigning the
# value, then restoring the cell style so no formatting is altered.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" "58.944.15"
Set-TextValue "E2" "  -1.19%  "
Set-TextValue "D3" "2.630.12"
Set-TextValue "E3" "  +1.09%  "
Set-TextValue "E4" "  +0.12%  "
Set-TextValue "D5" "513.74"
Set-TextValue "E5" "  +0.10%  "
Set-TextValue "D6" "144.46"
Set-TextValue "E6" "  -1.40%  "
Set-TextValue "D7" "0.996"
Set-TextValue "E7" "  -0.33%  "
Set-TextValue "D8" "0.571"
Set-TextValue "E8" "  +1.50%  "
Set-TextValue "D9" "2.656.68"
Set-TextValue "E9" "  +1.99%  "
Set-TextValue "D10" "6.32"
Set-TextValue "E10" "  +1.36%  "
Set-TextValue "E11" "  +2.18%  "
Set-TextValue "E12" "  +0.66%  "
Set-TextValue "E13" "  -1.42%  "
Set-TextValue "D14" "3.095.18"
Set-TextValue "E14" "  +1.24%  "
Set-TextValue "D15" "58.948.26"
Set-TextValue "E15" "  -1.18%  "
Set-TextValue "D16" "21.08"
Set-TextValue "E16" "  +0.78%  "
Set-TextValue "E17" "  +0.88%  "
Set-TextValue "D18" "2.654.01"
Set-TextValue "E18" "  +2.07%  "
Set-TextValue "E19" "  -0.41%  "
Set-TextValue "D20" "344.17"
Set-TextValue "E20" "  +2.11%  "
Set-TextValue "E21" "  +1.59%  "
Set-TextValue "D22" "6.09"
Set-TextValue "E22" "  +1.22%  "
Set-TextValue "D23" "0.999"
Set-TextValue "E23" "  +0.06%  "
Set-TextValue "D24" "60.94"
Set-TextValue "E24" "  +0.57%  "
Set-TextValue "D25" "0.420"
Set-TextValue "E25" "  +2.01%  "
Set-TextValue "D26" "2.759.52"
Set-TextValue "E26" "  +1.28%  "
Set-TextValue "D27" "0.995"
Set-TextValue "E27" "  -0.57%  "
Set-TextValue "D28" "0.160"
Set-TextValue "E28" "  +2.52%  "
Set-TextValue "D29" "0.0₃0807"
Set-TextValue "E29" "  +2.11%  "
Set-TextValue "E30" "  +2.77%  "
Set-TextValue "D31" "0.997"
Set-TextValue "E31" "  -0.28%  "
Set-TextValue "D32" "6.41"
Set-TextValue "E32" "  +8.93%  "
Set-TextValue "E33" "  +0.84%  "
Set-TextValue "D34" "18.91"
Set-TextValue "E34" "  +1.27%  "
Set-TextValue "D35" "148.91"
Set-TextValue "E35" "  -0.81%  "
Set-TextValue "E36" "  +13.30%  "
Set-TextValue "D37" "4.05"
Set-TextValue "E37" "  +4.25%  "
Set-TextValue "E38" "  +2.96%  "
Set-TextValue "D39" "0.855"
Set-TextValue "E39" "  +2.21%  "
Set-TextValue "D40" "36.49"
Set-TextValue "E40" "  -0.06%  "
Set-TextValue "E41" "  +3.62%  "
Set-TextValue "E42" "  +0.40%  "
Set-TextValue "D43" "280.93"
Set-TextValue "E43" "  -0.89%  "
Set-TextValue "D44" "0.615"
Set-TextValue "E44" "  -0.58%  "
Set-TextValue "E45" "  -0.53%  "
Set-TextValue "D46" "0.0985"
Set-TextValue "E46" "  -0.55%  "
Set-TextValue "D47" "19.49"
Set-TextValue "E47" "  +3.26%  "
Set-TextValue "D48" "0.0535"
Set-TextValue "E48" "  -0.67%  "
Set-TextValue "E49" "  -0.92%  "
Set-TextValue "E50" "  -0.07%  "
Set-TextValue "D51" "1.981.36"
Set-TextValue "E51" "  +2.74%  "
